# "Them cot nhiem vu cho tab gtn" - add a "task/role" column (D) next to the
# existing roster (A:C), giving each member a presentation-tab assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D values (order below matches the order the new shared strings
# were first introduced: Tim nghiem, Noi suy, Hoi quy + Gioi thieu nhom,
# Dao ham, Tich phan - row 5 reuses the same text as row 3).
$ws.Range("D2").Value = "Tab Tìm nghiệm"
$ws.Range("D1").Value = "Tab Nội suy"
$ws.Range("D3").Value = "Tab Hồi quy + Giới thiệu nhóm"
$ws.Range("D4").Value = "Tab Đạo hàm"
$ws.Range("D5").Value = "Tab Hồi quy + Giới thiệu nhóm"
$ws.Range("D6").Value = "Tab Tích phân"

# Widen columns B/C and size the new column D so the Vietnamese text fits.
$ws.Columns.Item(2).ColumnWidth = 15.0
$ws.Columns.Item(3).ColumnWidth = 20.333333333333332
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668

# Move the active selection off the data block (was A7:U26) onto F7.
$null = $ws.Range("F7").Select()
